# Insert a new weekly record row above current row 87 ("Hortaliza, Terminal La
# Palmera de La Serena - Jengibre"). This shifts all existing rows 87:141 down
# to 88:142 (Excel default insert behavior, which also extends the used range
# / dimension to A1:R142 and carries the date-column number format down with
# the shifted rows), and the freshly inserted row 87 is then populated with
# the new weekly observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 87:141 down to 88:142, inheriting formatting from the row above
# (this matches Excel's native "Insert Sheet Rows" behavior).
$ws.Rows.Item(87).Insert()

# Populate the newly inserted row 87 with the new record. Values mirror the
# former row 87 (now row 88) except for the Fecha (column D) and Volumen
# (column J), which carry the new observation's data.
$ws.Cells.Item(87, 1).Value = 8
$ws.Cells.Item(87, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(87, 3).Value = "Coquimbo"
$ws.Cells.Item(87, 4).Value = (Get-Date -Year 2023 -Month 6 -Day 19 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(87, 5).Value = 4
$ws.Cells.Item(87, 6).Value = 100114007
$ws.Cells.Item(87, 7).Value = "Jengibre"
$ws.Cells.Item(87, 8).Value = "Sin especificar"
$ws.Cells.Item(87, 9).Value = "Primera"
$ws.Cells.Item(87, 10).Value = 300
$ws.Cells.Item(87, 11).Value = 17000
$ws.Cells.Item(87, 12).Value = 18000
$ws.Cells.Item(87, 13).Value = 17500
$ws.Cells.Item(87, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(87, 15).Value = "Perú"
$ws.Cells.Item(87, 16).Value = 1346
$ws.Cells.Item(87, 17).Value = 13
$ws.Cells.Item(87, 18).Value = "Hortaliza"
